# Update the cryptocurrency "Price" (D) and "Volume(1h)" (E) columns with
# refreshed figures from the data source (GitHub Actions symbol-list sync).
#
# These columns hold plain text ("303.76", "4.08%", ...) rather than real
# numbers/percentages, so each cell is forced to Text format before the new
# value is written (otherwise Excel would auto-convert "303.52" to a number
# or "3.83%" to a percentage) and the style is reset back to Normal
# afterwards so no stray number-format style sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "303.52"
Set-TextValue $ws.Range("E2") "3.83%"
Set-TextValue $ws.Range("D3") "32.09"
Set-TextValue $ws.Range("E3") "8.47%"
Set-TextValue $ws.Range("D4") "5.240"
Set-TextValue $ws.Range("E4") "-0.74%"
Set-TextValue $ws.Range("D5") "0.07541"
Set-TextValue $ws.Range("E5") "5.18%"
Set-TextValue $ws.Range("D6") "7.936"
Set-TextValue $ws.Range("E6") "5.30%"
Set-TextValue $ws.Range("D7") "3.816"
Set-TextValue $ws.Range("E7") "6.44%"
Set-TextValue $ws.Range("D8") "1.514"
Set-TextValue $ws.Range("E8") "8.11%"
Set-TextValue $ws.Range("D9") "0.9229"
Set-TextValue $ws.Range("E9") "1.60%"
Set-TextValue $ws.Range("D10") "0.1701"
Set-TextValue $ws.Range("E10") "4.55%"
Set-TextValue $ws.Range("E11") "3.45%"
Set-TextValue $ws.Range("D12") "0.08034"
Set-TextValue $ws.Range("E12") "3.71%"
Set-TextValue $ws.Range("D13") "0.03037"
Set-TextValue $ws.Range("E13") "4.36%"
Set-TextValue $ws.Range("D14") "0.09890"
Set-TextValue $ws.Range("E14") "9.98%"
Set-TextValue $ws.Range("D15") "0.001505"
Set-TextValue $ws.Range("E15") "-5.25%"
Set-TextValue $ws.Range("D16") "0.04604"
Set-TextValue $ws.Range("E16") "1.80%"
Set-TextValue $ws.Range("D17") "0.006363"
Set-TextValue $ws.Range("E17") "2.34%"
Set-TextValue $ws.Range("D18") "3.465"
Set-TextValue $ws.Range("E18") "-0.33%"
Set-TextValue $ws.Range("D19") "2.232"
Set-TextValue $ws.Range("E19") "-0.02%"
Set-TextValue $ws.Range("E20") "1.47%"
Set-TextValue $ws.Range("D21") "0.1326"
Set-TextValue $ws.Range("E21") "-2.87%"
Set-TextValue $ws.Range("D22") "4.482"
Set-TextValue $ws.Range("E22") "11.00%"
Set-TextValue $ws.Range("E23") "1.77%"
Set-TextValue $ws.Range("D24") "0.001216"
Set-TextValue $ws.Range("E24") "0.98%"
Set-TextValue $ws.Range("D25") "0.004461"
Set-TextValue $ws.Range("E25") "5.00%"
Set-TextValue $ws.Range("D26") "0.0001399"
Set-TextValue $ws.Range("E26") "20.21%"
Set-TextValue $ws.Range("D27") "0.0001782"
Set-TextValue $ws.Range("E27") "6.13%"
Set-TextValue $ws.Range("D39") "0.01699"
Set-TextValue $ws.Range("E39") "2,502.27%"
Set-TextValue $ws.Range("D40") "0.04489"
Set-TextValue $ws.Range("E40") "1.06%"
Set-TextValue $ws.Range("D41") "0.006946"
Set-TextValue $ws.Range("E41") "-0.77%"
Set-TextValue $ws.Range("D42") "0.1353"
Set-TextValue $ws.Range("E42") "6.10%"
Set-TextValue $ws.Range("D43") "0.002078"
Set-TextValue $ws.Range("E43") "-5.44%"
Set-TextValue $ws.Range("D44") "0.01292"
Set-TextValue $ws.Range("E44") "-2.12%"
Set-TextValue $ws.Range("D45") "0.00006163"
Set-TextValue $ws.Range("E45") "5.96%"
Set-TextValue $ws.Range("D47") "0.01298"
Set-TextValue $ws.Range("E47") "0.46%"
